# Updated cryptos list values per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.448.13'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.037.57'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.74'
$ws.Range("E5").Value = '  +1.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.94'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.031.42'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.148'
$ws.Range("E10").Value = '  -3.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.07'
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.449'
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000222'
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.61'
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.535.32'
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.461.49'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.039.67'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.109'
$ws.Range("E18").Value = '  -2.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.66'
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '481.42'
$ws.Range("E20").Value = '  +3.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.28'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.676'
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.07'
$ws.Range("E23").Value = '  +2.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.82'
$ws.Range("E24").Value = '  +4.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.16'
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.72'
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.78'
$ws.Range("E28").Value = '  +0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +5.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.85'
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.69'
$ws.Range("E33").Value = '  +6.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.37'
$ws.Range("E34").Value = '  +4.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.06'
$ws.Range("E35").Value = '  -6.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.89'
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '453.90'
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.179.20'
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0803'
$ws.Range("E39").Value = '  +2.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0386'
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.118'
$ws.Range("E41").Value = '  +2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.13'
$ws.Range("E42").Value = '  +1.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.46'
$ws.Range("E43").Value = '  -2.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.49'
$ws.Range("E44").Value = '  +7.44%  '
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.245'
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.109'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.98'
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '115.87'
$ws.Range("E49").Value = '  -4.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₃0498'
$ws.Range("E50").Value = '  -2.33%  '
$ws.Range("E51").Value = '  +3.76%  '
